# Adding the front and back heading to the card for more flexibility
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: statesAndCapitals
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("statesAndCapitals")

# D1: standalone formula gets "State:" / "Capital:" headings prefixed.
$ws1.Range("D1").Formula = '=CONCATENATE("new card(""State:"",""", A1, """,""Capital:"",""",B1, """),")'

# D2:D50: shared formula (si=0) gets the same headings.
$ws1.Range("D2:D50").Formula = '=CONCATENATE("new card(""State:"",""", A2, """,""Capital:"",""",B2, """),")'

# Column D widens to fit the new, longer text.
$ws1.Columns.Item(4).ColumnWidth = 50

# View: scrolled down and D1:D50 selected.
$ws1.Activate()
$ws1.Range("D1:D50").Select()
$win1 = $excel.ActiveWindow
$win1.ScrollRow = 21
$win1.ScrollColumn = 1

# ---------------------------------------------------------------------------
# Sheet: multiplication
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("multiplication")

# E1: standalone formula gets "Question:" / "Answer:" (with colon) headings.
$ws2.Range("E1").Formula = '=CONCATENATE("new card(""Question:"",""", A1, " X ", B1, """,""Answer:"",""",C1, """),")'

# E2:E65: shared formula (si=1) gets "Question" / "Answer" (no colon).
$ws2.Range("E2:E65").Formula = '=CONCATENATE("new card(""Question"",""", A2, " X ", B2, """,""Answer"",""",C2, """),")'

# E66:E129: shared formula (si=3) gets "Question" / "Answer" (no colon).
$ws2.Range("E66:E129").Formula = '=CONCATENATE("new card(""Question"",""", A66, " X ", B66, """,""Answer"",""",C66, """),")'

# E130:E144: shared formula (si=5) gets "Question" / "Answer" (no colon).
$ws2.Range("E130:E144").Formula = '=CONCATENATE("new card(""Question"",""", A130, " X ", B130, """,""Answer"",""",C130, """),")'

# View: scrolled down further and G120 selected.
$ws2.Activate()
$ws2.Range("G120").Select()
$win2 = $excel.ActiveWindow
$win2.ScrollRow = 115
$win2.ScrollColumn = 1
